$d = $word.ActiveDocument

# 1. Wrap the first bullet's text in a bookmark (as seen in the target
#    OOXML: bookmarkStart right before the run in paragraph 1, bookmarkEnd
#    at the paragraph boundary right after it).
$p1 = $d.Paragraphs(1)
$d.Bookmarks.Add("_Hlk176169549", $p1.Range) | Out-Null

# 2. Trim the trailing space in ". Timestamps fx. " -> ". Timestamps fx."
$d.Content.Find.Execute(". Timestamps fx. ", $true, $false, $false, $false, $false, $true, 1, $false, ". Timestamps fx.", 2) | Out-Null

# 3. Split the bullet so "Coregistrering af tryk-plots ..." becomes its own
#    list item, matching the new <w:p> that starts right before it.
$sel = $word.Selection
$sel.Find.Execute("Coregistrering af tryk-plots", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$sel.Collapse(1)
$sel.InsertParagraphBefore()
